$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.117.89'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '1.646.14'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").Value = '''216.88'
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").Value = '''0.5151'
$ws.Range("E6").Value = '  +1.74%  '
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("D8").Value = '''0.2605'
$ws.Range("E8").Value = '  +1.14%  '
$ws.Range("D9").Value = '''0.06409'
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").Value = '''19.97'
$ws.Range("E10").Value = '  +1.06%  '
$ws.Range("D11").Value = '''0.07789'
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.676.48'
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.325'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").Value = '''0.5529'
$ws.Range("E14").Value = '  +1.32%  '
$ws.Range("D15").Value = '''65.04'
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").Value = '26.125.13'
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("D18").Value = '''1.003'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("D19").Value = '''201.27'
$ws.Range("E19").Value = '  +1.62%  '
$ws.Range("D20").Value = '''4.505'
$ws.Range("E20").Value = '  +2.09%  '
$ws.Range("D21").Value = '''10.05'
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").Value = '''6.149'
$ws.Range("E22").Value = '  +1.68%  '
$ws.Range("D23").Value = '''1.004'
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").Value = '''1.904'
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("D25").Value = '''142.61'
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("D26").Value = '''0.1223'
$ws.Range("E26").Value = '  +6.48%  '
$ws.Range("D27").Value = '''6.939'
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").Value = '''15.78'
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").Value = '''0.04951'
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("D31").Value = '''3.343'
$ws.Range("E31").Value = '  +2.15%  '
$ws.Range("D32").Value = '''3.275'
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("D33").Value = '''1.556'
$ws.Range("E33").Value = '  +0.90%  '
$ws.Range("D34").Value = '''2.380'
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").Value = '''0.9268'
$ws.Range("E35").Value = '  +3.58%  '
$ws.Range("D36").Value = '''0.5636'
$ws.Range("E36").Value = '  +1.68%  '
$ws.Range("D37").Value = '''2.593'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").Value = '1.119.47'
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("D39").Value = '''0.01581'
$ws.Range("E39").Value = '  +1.13%  '
$ws.Range("D40").Value = '''1.002'
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").Value = '''2.542'
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("D42").Value = '''5.615'
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("D43").Value = '''0.8147'
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").Value = '''100.08'
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").Value = '0.0₈124'
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("D46").Value = '1.783.66'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("D47").Value = '''0.4546'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").Value = '''55.69'
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("D49").Value = '''1.004'
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = '''0.05279'
$ws.Range("E50").Value = '  +3.66%  '
$ws.Range("D51").Value = '''0.09658'
$ws.Range("E51").Value = '  +0.93%  '
